$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 94582
$ws.Range("B2").Value = "Gustavo Campos"
$ws.Range("C2").Value = "Engenharia"
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 45096
$ws.Range("G2").Value = 12447.97

$ws.Range("A3").Value = 49998
$ws.Range("B3").Value = "Alexandre Costa"
$ws.Range("C3").Value = "Recursos Humanos"
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 45098
$ws.Range("G3").Value = 6808.09

$ws.Range("A4").Value = 57452
$ws.Range("B4").Value = "Manuela Fogaça"
$ws.Range("C4").Value = "Financeiro"
$ws.Range("D4").Value = "Viagem de negócios"
$ws.Range("E4").Value = 7
$ws.Range("F4").Value = 45079
$ws.Range("G4").Value = 12059.61

$ws.Range("A5").Value = 69626
$ws.Range("B5").Value = "João Vitor Fernandes"
$ws.Range("C5").Value = "Recursos Humanos"
$ws.Range("D5").Value = "Viagem de negócios"
$ws.Range("E5").Value = 6
$ws.Range("F5").Value = 45078
$ws.Range("G5").Value = 10971.33

$ws.Range("A6").Value = 69539
$ws.Range("B6").Value = "Ana Beatriz Nascimento"
$ws.Range("C6").Value = "Engenharia"
$ws.Range("D6").Value = "Doença"
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 45091
$ws.Range("G6").Value = 5843.19

$ws.Range("A7").Value = 98168
$ws.Range("B7").Value = "Lucas Cardoso"
$ws.Range("C7").Value = "Financeiro"
$ws.Range("D7").Value = "Viagem de negócios"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 45097
$ws.Range("G7").Value = 8161.13

$ws.Range("A8").Value = 768
$ws.Range("B8").Value = "Daniel Fogaça"
$ws.Range("C8").Value = "Atendimento ao Cliente"
$ws.Range("D8").Value = "Problemas pessoais"
$ws.Range("F8").Value = 45085
$ws.Range("G8").Value = 9141.34

$ws.Range("A9").Value = 36287
$ws.Range("B9").Value = "Srta. Maysa Barros"
$ws.Range("C9").Value = "Vendas"
$ws.Range("D9").Value = "Outros"
$ws.Range("E9").Value = 7
$ws.Range("F9").Value = 45078
$ws.Range("G9").Value = 5553.8

$ws.Range("A10").Value = 37129
$ws.Range("B10").Value = "Sra. Valentina Santos"
$ws.Range("C10").Value = "Vendas"
$ws.Range("D10").Value = "Doença"
$ws.Range("E10").Value = 6
$ws.Range("F10").Value = 45084
$ws.Range("G10").Value = 4433.91

$ws.Range("A11").Value = 87742
$ws.Range("B11").Value = "João Lucas Sales"
$ws.Range("C11").Value = "Operações"
$ws.Range("D11").Value = "Doença"
$ws.Range("F11").Value = 45081
$ws.Range("G11").Value = 8973.389999999999
